$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New D-column values for rows 128-132
$ws.Range("D128").Value = 0.5614898640156477
$ws.Range("D129").Value = 0.7804486233837438
$ws.Range("D130").Value = 0.7918746591162638
$ws.Range("D131").Value = 0.7115302101128926
$ws.Range("D132").Value = 0.7732527029753358

# New C-column values for rows 133-139
$ws.Range("C133").Value = 0.5947585845621888
$ws.Range("C134").Value = -0.1471494033888783
$ws.Range("C135").Value = 0.2315426859737066
$ws.Range("C136").Value = 0.5186180303197108
$ws.Range("C137").Value = 0.5462623555331837
$ws.Range("C138").Value = 0.3054124296933831
$ws.Range("C139").Value = 0.3912781359200057

# New B-column value for row 140
$ws.Range("B140").Value = 0.09995295491779371
